$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 43) with the latest Adafruit IO reading.
$row = 43

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"

# Column C holds a numeric-looking reading ("25") that must be stored as
# text (matching every other row in this sheet). Prefix with an
# apostrophe so Excel treats it as text instead of auto-converting it to
# a number, then clear the resulting "quote prefix" formatting so no
# stray style is left behind on the cell.
$ws.Cells.Item($row, 3).Value = "'25"
$ws.Cells.Item($row, 3).ClearFormats()

$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"
